$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2517.5881
$ws.Range("I32").Value = 1676.375
$ws.Range("J32").Value = 3265.3333
$ws.Range("K32").Value = 1676.375
$ws.Range("L32").Value = 3265.3333
$ws.Range("M32").Value = -1350.375
$ws.Range("N32").Value = -3917.3333
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H76").Value = 4000
$ws.Range("I76").Value = 4000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 4000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3685
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 4000
$ws.Range("I79").Value = 4000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 4000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2908
$ws.Range("N79").ClearContents()
$ws.Range("H98").Value = 1756.4375
$ws.Range("I98").Value = 1756.4375
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1756.4375
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -258.4375
$ws.Range("N98").ClearContents()
$ws.Range("H113").Value = 2377.7368
$ws.Range("I113").Value = 1823.625
$ws.Range("J113").Value = 5333
$ws.Range("K113").Value = 1823.625
$ws.Range("L113").Value = 5333
$ws.Range("M113").Value = 1430.375
$ws.Range("N113").Value = -11841
$ws.Range("H122").Value = 1756.4375
$ws.Range("I122").Value = 1756.4375
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5269.3125
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2819.3125
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 5714.4546
$ws.Range("I132").Value = 5982.75
$ws.Range("K132").Value = 17948.25
$ws.Range("M132").Value = -15418.25
$ws.Range("H138").Value = 1514.1
$ws.Range("I138").Value = 1608.4546
$ws.Range("J138").Value = 1398.7778
$ws.Range("K138").Value = 4825.3638
$ws.Range("L138").Value = 4196.3334
$ws.Range("M138").Value = 314.6361999999999
$ws.Range("N138").Value = -14476.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 290.14285
$ws.Range("I5").Value = 327
$ws.Range("K5").Value = 327
$ws.Range("M5").Value = -215
$ws.Range("H32").Value = 7352.276
$ws.Range("I32").Value = 6600.6294
$ws.Range("K32").Value = 6600.6294
$ws.Range("M32").Value = -6313.6294
$ws.Range("H39").Value = 11000
$ws.Range("I39").Value = 11000
$ws.Range("K39").Value = 11000
$ws.Range("M39").Value = -10480
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 3297.6667
$ws.Range("I97").Value = 1957.2
$ws.Range("K97").Value = 1957.2
$ws.Range("M97").Value = -1461.2
$ws.Range("H102").Value = 1499.6
$ws.Range("I102").Value = 1499.6
$ws.Range("K102").Value = 1499.6
$ws.Range("M102").Value = 122.4000000000001
$ws.Range("H132").Value = 2567.4285
$ws.Range("I132").Value = 2630.8
$ws.Range("K132").Value = 7892.400000000001
$ws.Range("M132").Value = -5362.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 290.14285
$ws.Range("I4").Value = 327
$ws.Range("K4").Value = 327
$ws.Range("M4").Value = -212
$ws.Range("H20").Value = 3195.2307
$ws.Range("I20").Value = 3379.4
$ws.Range("K20").Value = 3379.4
$ws.Range("M20").Value = -3132.4
$ws.Range("H94").Value = 1539.3077
$ws.Range("I94").Value = 2854.6
$ws.Range("J94").Value = 1226.1428
$ws.Range("K94").Value = 2854.6
$ws.Range("L94").Value = 1226.1428
$ws.Range("M94").Value = -2403.6
$ws.Range("N94").Value = -2128.1428
$ws.Range("H99").Value = 2772.875
$ws.Range("J99").Value = 2997.5715
$ws.Range("L99").Value = 2997.5715
$ws.Range("N99").Value = -5993.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2825.7222
$ws.Range("I31").Value = 2112
$ws.Range("J31").Value = 3717.875
$ws.Range("K31").Value = 2112
$ws.Range("L31").Value = 3717.875
$ws.Range("M31").Value = -1817
$ws.Range("N31").Value = -4307.875
$ws.Range("H34").Value = 2825.7222
$ws.Range("I34").Value = 2112
$ws.Range("J34").Value = 3717.875
$ws.Range("K34").Value = 2112
$ws.Range("L34").Value = 3717.875
$ws.Range("M34").Value = -1910
$ws.Range("N34").Value = -4121.875
$ws.Range("H99").Value = 6273.8
$ws.Range("I99").Value = 6273.8
$ws.Range("K99").Value = 6273.8
$ws.Range("M99").Value = -4775.8
$ws.Range("H121").Value = 19332.666
$ws.Range("J121").Value = 19332.666
$ws.Range("L121").Value = 19332.666
$ws.Range("N121").Value = -21952.666
$ws.Range("H122").Value = 1962.6666
$ws.Range("I122").Value = 1962.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5887.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3437.9998
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 6273.8
$ws.Range("I126").Value = 6273.8
$ws.Range("K126").Value = 18821.4
$ws.Range("M126").Value = -16351.4
$ws.Range("H132").Value = 2793.6
$ws.Range("I132").Value = 2398.25
$ws.Range("K132").Value = 7194.75
$ws.Range("M132").Value = -4664.75
$ws.Range("H134").Value = 2460
$ws.Range("I134").Value = 2165.375
$ws.Range("K134").Value = 6496.125
$ws.Range("M134").Value = -3961.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 82.76470999999999
$ws.Range("I2").Value = 87.23077000000001
$ws.Range("J2").Value = 68.25
$ws.Range("K2").Value = 523.38462
$ws.Range("L2").Value = 409.5
$ws.Range("M2").Value = -410.38462
$ws.Range("N2").Value = -635.5
$ws.Range("H112").Value = 10000
$ws.Range("J112").Value = 10000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32216
$ws.Range("H137").Value = 2500
$ws.Range("J137").Value = 2500
$ws.Range("L137").Value = 7500
$ws.Range("N137").Value = -17700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 17627500
$ws.Range("J11").Value = 3673334.8
$ws.Range("L11").Value = 3673334.8
$ws.Range("N11").Value = -3673612.8
$ws.Range("H70").Value = 9499.6
$ws.Range("I70").Value = 9249
$ws.Range("K70").Value = 9249
$ws.Range("M70").Value = -8979
$ws.Range("H73").Value = 9499.6
$ws.Range("I73").Value = 9249
$ws.Range("K73").Value = 9249
$ws.Range("M73").Value = -8313
$ws.Range("H132").Value = 3823
$ws.Range("I132").Value = 3787.6
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 11362.8
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -8832.799999999999
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1999.5
$ws.Range("I7").Value = 1999.5
$ws.Range("K7").Value = 1999.5
$ws.Range("M7").Value = -1887.5
$ws.Range("H46").Value = 1446.75
$ws.Range("I46").Value = 1006.6429
$ws.Range("J46").Value = 4527.5
$ws.Range("K46").Value = 1006.6429
$ws.Range("L46").Value = 4527.5
$ws.Range("M46").Value = -818.6429000000001
$ws.Range("N46").Value = -4903.5
$ws.Range("H61").Value = 1491.8572
$ws.Range("I61").Value = 1582.3334
$ws.Range("J61").Value = 949
$ws.Range("K61").Value = 1582.3334
$ws.Range("L61").Value = 949
$ws.Range("M61").Value = -1380.3334
$ws.Range("N61").Value = -1353
$ws.Range("H113").Value = 1491.8572
$ws.Range("I113").Value = 1582.3334
$ws.Range("J113").Value = 949
$ws.Range("K113").Value = 1582.3334
$ws.Range("L113").Value = 949
$ws.Range("M113").Value = 587.6666
$ws.Range("N113").Value = -5289
$ws.Range("H122").Value = 3660.5
$ws.Range("I122").Value = 3592.6
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 10777.8
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8327.799999999999
$ws.Range("N122").Value = -16900
$ws.Range("H126").Value = 1999.5
$ws.Range("I126").Value = 1999.5
$ws.Range("K126").Value = 5998.5
$ws.Range("M126").Value = -3528.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 346666
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 346666
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 346666
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -346890
$ws.Range("H61").Value = 500000
$ws.Range("I61").Value = 500000
$ws.Range("K61").Value = 500000
$ws.Range("M61").Value = -499708
$ws.Range("H132").Value = 3164.6155
$ws.Range("I132").Value = 2595
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 7785
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -5255
$ws.Range("N132").Value = -35060
